$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# Update Version value
$ws1.Range("B3").Value = "6.0.0"

# Update Date value
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value (was empty)
$ws1.Range("B9").Value = "Alvearie Team"

# Replace the duplicate "Contact" row (row 11) with "Jurisdiction" / "United States of America",
# and shift remaining rows (Description, Purpose, Copyright, Immutable) up by one.
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

$ws1.Range("A11").Value = "Description"
$ws1.Range("B11").Value = "Codes indicating methods used to identify matched entities"

$ws1.Range("A12").Value = "Purpose"
$ws1.Range("B12").Value = $null

$ws1.Range("A13").Value = "Copyright"
$ws1.Range("B13").Value = $null

$ws1.Range("A14").Value = "Immutable"
$ws1.Range("B14").Value = "BooleanType[null]"

# Delete now-empty row 15 (the old last row content moved to row14 already, so delete the leftover row 15)
$ws1.Range("A15:B15").Delete()
